$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic rotation of the data rows 4-6:
#   new row 4 <- old row 6
#   new row 5 <- old row 4
#   new row 6 <- old row 5
# Coordinates (Q/R) are also rounded to whole numbers, and the
# Starttid/Sluttid ("00:00") cells (Z/AB) are removed for every row.

# ---- Row 4 (becomes former row 6's data: Blåsippa) ----
$ws.Range("A4").Value = 111782567
$ws.Range("B4").Value = 98535
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 222498
$ws.Range("F4").Value = "Blåsippa"
$ws.Range("G4").Value = "Hepatica nobilis"
$ws.Range("H4").Value = "Schreb."
$ws.Range("I4").Value = "'5"
$ws.Range("J4").Value = "plantor/tuvor"
$ws.Range("Q4").Value = 573909
$ws.Range("R4").Value = 6303235
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()

# ---- Row 5 (becomes former row 4's data: Skogsalm) ----
$ws.Range("A5").Value = 111782565
$ws.Range("B5").Value = 100532
$ws.Range("D5").Value = "CR"
$ws.Range("E5").Value = 223246
$ws.Range("F5").Value = "Skogsalm"
$ws.Range("G5").Value = "Ulmus glabra"
$ws.Range("H5").Value = "Huds."
# I5 stays "1" (unchanged by the edit), so it is left untouched.
$ws.Range("J5").Value = "plantor/tuvor"
$ws.Range("Q5").Value = 573877
$ws.Range("R5").Value = 6303226
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()

# ---- Row 6 (becomes former row 5's data: Myskmadra) ----
$ws.Range("A6").Value = 111782566
$ws.Range("B6").Value = 103369
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 221423
$ws.Range("F6").Value = "Myskmadra"
$ws.Range("G6").Value = "Galium odoratum"
$ws.Range("H6").Value = "(L.) Scop."
$ws.Range("I6").Value = "'1"
$ws.Range("J6").Value = "m²"
$ws.Range("Q6").Value = 573877
$ws.Range("R6").Value = 6303234
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()
